$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("2.S1 Back Savers")

# Update the "units produced" inputs (C11:D11) that feed the SUMPRODUCT formulas
$ws.Range("C11").Value = 1000
$ws.Range("D11").Value = 974.99999995125029

# Recalculate so dependent formulas (E7, E8, G11) pick up new values
$excel.CalculateFullRebuild()

# Move the active selection to G11 as in the edited workbook
$ws.Activate()
$ws.Range("G11").Select()
